$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (row 9) down onto the
# three new rows, one column at a time so each column keeps its own
# (slightly different) cell style, then layer Bold on top of the copied
# formatting to get the new "bold" style variants used for the latest
# entries.
$ws.Range("A9").Copy($ws.Range("A10"))
$ws.Range("B9").Copy($ws.Range("B10"))
$ws.Range("C9").Copy($ws.Range("C10"))

$ws.Range("A9").Copy($ws.Range("A11"))
$ws.Range("B9").Copy($ws.Range("B11"))
$ws.Range("C9").Copy($ws.Range("C11"))

$ws.Range("A9").Copy($ws.Range("A12"))
$ws.Range("B9").Copy($ws.Range("B12"))
$ws.Range("C9").Copy($ws.Range("C12"))

$ws.Range("A10:C12").Font.Bold = $true

# New log entries
$ws.Range("A10").Value = 45695
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = "Coming up with the possible endpoints"

$ws.Range("A11").Value = 45696
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = "Setting up the backend and writing the project report"

$ws.Range("A12").Value = 45696
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = "Committing to github repo"
